# Recreate the mothertongues "Skip" worksheet: a second CSV-imported table
# (header row + data rows) alongside Sheet1's existing query-table import,
# with a matching "data_1" defined name pointing at the new range.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item("Sheet1")

# Add the new sheet right after Sheet1 and name it "Skip"
$skip = $wb.Worksheets.Add($null, $sheet1)
$skip.Name = "Skip"

# Header row
$skip.Range("A1").Value = "ID"
$skip.Range("B1").Value = "definition"
$skip.Range("C1").Value = "word"
$skip.Range("D1").Value = "audio_desc"
$skip.Range("E1").Value = "audio_fn"
$skip.Range("F1").Value = "theme"
$skip.Range("G1").Value = "secondary_theme"

# Data rows (same records as Sheet1, offset by the header + re-numbered ID column)
$skip.Range("A2").Value = 1
$skip.Range("B2").Value = "tree"
$skip.Range("C2").Value = "træ"
$skip.Range("D2").Value = "Aidan Pine"
$skip.Range("E2").Value = "tree.mp3"
$skip.Range("F2").Value = "plants"
$skip.Range("G2").Value = "noun"

$skip.Range("A3").Value = 2
$skip.Range("B3").Value = "word"
$skip.Range("C3").Value = "ord"
$skip.Range("D3").Value = "Aidan Pine"
$skip.Range("E3").Value = "ord.mp3"
$skip.Range("F3").Value = "abstract"
$skip.Range("G3").Value = "noun"

$skip.Range("A4").Value = 3
$skip.Range("B4").Value = "hello"
$skip.Range("C4").Value = "hej"
$skip.Range("D4").Value = "Aidan Pine"
$skip.Range("E4").Value = "hej.mp3"
$skip.Range("F4").Value = "greetings"
$skip.Range("G4").Value = "interjection"

$skip.Range("A5").Value = 4
$skip.Range("B5").Value = "goodbye"
$skip.Range("C5").Value = "farvel"
$skip.Range("F5").Value = "greetings"
$skip.Range("G5").Value = "interjection"

# Defined name "data_1" (localSheetId scoped to Skip) mirroring "data" on Sheet1
$skip.Names.Add("data_1", "=Skip!`$A`$2:`$G`$5")

# Leave selections where the source workbook had them: Skip on B8, Sheet1 back
# on top as the active sheet with D16 selected.
$skip.Range("B8").Select()

$sheet1.Activate()
$sheet1.Range("D16").Select()
